$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target cluster labels based on the revised NATMI run (Dr Hou advice)
$sCs = "sCs"
$Cd5l = "Cd5l"
$Cd5 = "Cd5"
$ECs = "ECs"
$FAPs = "FAPs"
$M2 = "M2"

function Set-Row {
    param($row, $a, $b, $c, $d, $vals)
    $ws.Range("A$row").Value = $a
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Value = $d
    $cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}

# Row 2: M2 -> ECs
Set-Row 2 $M2 $Cd5l $Cd5 $ECs @(3, 1, 3.311719666666667, 9.935159000000001, 0.9364811982846314, 0.9364811982846314, 3, 1, 1.439038333333333, 4.317115, 0.3449011277814652, 0.3449011277814653, 4.765691549587223, 42.891223946285, 0.3229934214345073, 0.3229934214345074)

# Row 3: M2 -> FAPs
Set-Row 3 $M2 $Cd5l $Cd5 $FAPs @(3, 1, 3.311719666666667, 9.935159000000001, 0.9364811982846314, 0.9364811982846314, 3, 1, 0.982904, 2.948712, 0.2355772533978687, 0.2355772533978687, 3.255102507245334, 29.295922565208, 0.2206136685506383, 0.2206136685506383)

# Row 4: M2 -> M2
Set-Row 4 $M2 $Cd5l $Cd5 $M2 @(3, 1, 3.311719666666667, 9.935159000000001, 0.9364811982846314, 0.9364811982846314, 3, 1, 1.621942333333333, 4.865826999999999, 0.3887385950778479, 0.3887385950778479, 5.371418323499222, 48.342764911493, 0.3640463853379871, 0.3640463853379871)

# Row 5: M2 -> sCs
Set-Row 5 $M2 $Cd5l $Cd5 $sCs @(3, 1, 3.311719666666667, 9.935159000000001, 0.9364811982846314, 0.9364811982846314, 3, 1, 0.1284366666666667, 0.38531, 0.03078302374281815, 0.03078302374281815, 0.4253462349211111, 3.82811611429, 0.0288277229614986, 0.0288277229614986)

# Row 6: sCs -> ECs
Set-Row 6 $sCs $Cd5l $Cd5 $ECs @(1, 0.3333333333333333, 0.2246243333333333, 0.6738730000000001, 0.06351880171536857, 0.06351880171536857, 3, 1, 1.439038333333333, 4.317115, 0.3449011277814652, 0.3449011277814653, 0.3232430262661111, 2.909187236395, 0.02190770634695789, 0.02190770634695789)

# Row 7: sCs -> FAPs
Set-Row 7 $sCs $Cd5l $Cd5 $FAPs @(1, 0.3333333333333333, 0.2246243333333333, 0.6738730000000001, 0.06351880171536857, 0.06351880171536857, 3, 1, 0.982904, 2.948712, 0.2355772533978687, 0.2355772533978687, 0.2207841557306667, 1.987057401576, 0.01496358484723036, 0.01496358484723036)

# Row 8: sCs -> M2
Set-Row 8 $sCs $Cd5l $Cd5 $M2 @(1, 0.3333333333333333, 0.2246243333333333, 0.6738730000000001, 0.06351880171536857, 0.06351880171536857, 3, 1, 1.621942333333333, 4.865826999999999, 0.3887385950778479, 0.3887385950778479, 0.3643277153301111, 3.278949437971, 0.02469220973986077, 0.02469220973986077)

# Row 9: sCs -> sCs
Set-Row 9 $sCs $Cd5l $Cd5 $sCs @(1, 0.3333333333333333, 0.2246243333333333, 0.6738730000000001, 0.06351880171536857, 0.06351880171536857, 3, 1, 0.1284366666666667, 0.38531, 0.03078302374281815, 0.03078302374281815, 0.02885000062555556, 0.25965000563, 0.001955300781319549, 0.001955300781319549)
